$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("GanttChart")
$ws.Range("I4").Value = 43
